# edit.ps1
# Applies the "rebuilt with updated templates" change described by the XML
# diff: inserts several new "Project"/"Ilias" subsections, a brand-new
# SC-12 and SC-15 control section, and fixes the "S3 buckets" wording.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Inserts a new, empty paragraph immediately before the paragraph that is
# currently at $beforeIdx (1-based), gives it the requested style, and
# returns the new paragraph object. After the call, the paragraph that used
# to be at $beforeIdx now sits at $beforeIdx + 1.
function New-ParaBefore($doc, $beforeIdx, $style) {
    $refPara = $doc.Paragraphs.Item($beforeIdx)
    $r = $doc.Range($refPara.Range.Start, $refPara.Range.Start)
    $r.InsertParagraphBefore()
    $newPara = $doc.Paragraphs.Item($beforeIdx)
    $newPara.Range.ParagraphFormat.Style = $style
    return $newPara
}

# Appends plain text to the end of the document range at $pos, optionally
# applying a character style (e.g. VerbatimChar) and/or bold, then returns
# the position right after the inserted text.
function Add-Run($doc, $pos, $text, $charStyle, $bold) {
    if ($text.Length -eq 0) {
        return $pos
    }
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
    $newPos = $pos + $text.Length
    $rs = $doc.Range($pos, $newPos)
    if ($charStyle) {
        $rs.Style = $charStyle
    }
    if ($bold) {
        $rs.Font.Bold = $true
    } else {
        $rs.Font.Bold = $false
    }
    return $newPos
}

# Appends a manual line break (renders as <w:br/>), returns new position.
function Add-LineBreak($doc, $pos) {
    $r = $doc.Range($pos, $pos)
    $r.InsertBreak(6)  # wdLineBreak
    return $pos + 1
}

# Wraps the text range [start, end) with a bookmark of the given name.
function Add-Bookmark($doc, $start, $end, $name) {
    $rng = $doc.Range($start, $end)
    $doc.Bookmarks.Add($name, $rng) | Out-Null
}

# Builds a Heading5 "title" paragraph (bookmarked) immediately before
# $beforeIdx, returns the index of the paragraph following the new one
# (i.e. the new value to use as the next insertion point).
function Add-Heading5($doc, $beforeIdx, $title, $bookmarkName) {
    $p = New-ParaBefore $doc $beforeIdx "Heading5"
    $s = $p.Range.Start
    $e = Add-Run $doc $s $title $null $false
    Add-Bookmark $doc $s $e $bookmarkName
    return $beforeIdx + 1
}

# Builds a Heading3 "title" paragraph (bookmarked) immediately before
# $beforeIdx, returns the next insertion index.
function Add-Heading3($doc, $beforeIdx, $title, $bookmarkName) {
    $p = New-ParaBefore $doc $beforeIdx "Heading3"
    $s = $p.Range.Start
    $e = Add-Run $doc $s $title $null $false
    Add-Bookmark $doc $s $e $bookmarkName
    return $beforeIdx + 1
}

# Builds a FirstParagraph made of one or more runs. $runs is an array of
# hashtables: @{ Text = "..."; Bold = $true/$false }
function Add-FirstParagraph($doc, $beforeIdx, $runs) {
    $p = New-ParaBefore $doc $beforeIdx "FirstParagraph"
    $pos = $p.Range.Start
    foreach ($run in $runs) {
        $pos = Add-Run $doc $pos $run.Text $null $run.Bold
    }
    return $beforeIdx + 1
}

# Builds a SourceCode paragraph made of one or more (text, isBreak) lines,
# each non-break chunk styled as VerbatimChar, mirroring the original
# document's "source code" blocks.
function Add-SourceCode($doc, $beforeIdx, $lines) {
    $p = New-ParaBefore $doc $beforeIdx "SourceCode"
    $pos = $p.Range.Start
    for ($i = 0; $i -lt $lines.Count; $i++) {
        if ($i -gt 0) {
            $pos = Add-LineBreak $doc $pos
        }
        $pos = Add-Run $doc $pos $lines[$i] "VerbatimChar" $false
    }
    return $beforeIdx + 1
}

# ---------------------------------------------------------------------
# 1) New "Project" subsection right before "SC-5: Denial Of Service
#    Protection" (after the existing CivicActions paragraph).
# ---------------------------------------------------------------------

$idx = 9   # "SC-5: Denial Of Service Protection" Heading3
$idx = Add-Heading5 $d $idx "Project" "project"
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "System and communications protection policy and procedures are formally documented in the None and the Project SSP. The Department reviews and updates the policy as necessary and has been continually updated since April 2008."; Bold = $false },
    @{ Text = " "; Bold = $false },
    @{ Text = "This is Agency common control. More data about implementation can be obtained from the Agency common control catalog."; Bold = $false }
)

# ---------------------------------------------------------------------
# 2) New "Ilias" + "Project" subsections right before "SC-7: Boundary
#    Protection" (after the existing Drupal paragraph).
# ---------------------------------------------------------------------

$idx = 16  # "SC-7: Boundary Protection" Heading3 (shifted by +2 from step 1)
$idx = Add-Heading5 $d $idx "Ilias" "ilias"
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "Ilias has a manual ability to block IP addresses in cases where attacks bypass cloud protection. This is managed by CivicActions Operations."; Bold = $false }
)
$idx = Add-Heading5 $d $idx "Project" "project-1"
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "The Project system is configured to reduce vulnerabilities in its operating system and applications to protect against Denial of Service (DoS) attacks."; Bold = $false },
    @{ Text = " "; Bold = $false },
    @{ Text = "The Project support staff ensures the system is protected against or limits the effect of DoS attacks as specified in the None."; Bold = $false }
)

# ---------------------------------------------------------------------
# 3) New "Ilias" + "Project" subsections right before the "a" Heading4
#    (after the existing Drupal/SELinux paragraph).
# ---------------------------------------------------------------------

$idx = 25  # "a" Heading4 (shifted by +4 from step 2)
$idx = Add-Heading5 $d $idx "Ilias" "ilias-1"
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "Ilias, when deployed on SELinux in full enforcing mode, minimizes the number of services and computing nodes that are exposed to the Internet. Ilias employs both the AWS platform safeguards and the Ilias logging in monitoring and recording system events. All other computing nodes used in the system are isolated within AWS."; Bold = $false }
)
$idx = Add-Heading5 $d $idx "Project" "project-2"
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "The Project system has monitored and controlled communications at the external boundary of the information system and at key internal boundaries within the system, where appropriate. The Project allocates publicly accessible information system components (e.g., public web servers) specific IP address and port combinations. Public access into the organization’s internal networks is prevented except as appropriately mediated."; Bold = $false }
)

# ---------------------------------------------------------------------
# 4) Wording fix: "S3 buckets" -> "Amazon S3 buckets" in the first AWS
#    paragraph under heading "a".
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "VPCs, subnets and S3 buckets",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "VPCs, subnets and Amazon S3 buckets", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Brand-new "SC-12: Cryptographic Key Establishment And Management"
#    section, right before "SC-13: Cryptographic Protection".
# ---------------------------------------------------------------------

$idx = 38  # "SC-13: Cryptographic Protection" Heading3 (shifted by +4 from step 3)
$idx = Add-Heading3 $d $idx "SC-12: Cryptographic Key Establishment And Management" "sc-12-cryptographic-key-establishment-and-management"
$idx = Add-SourceCode $d $idx @(
    "The organization establishes and manages cryptographic keys for required cryptography employed within the information system in accordance with [Assignment: organization-defined requirements for key generation, distribution, storage, access, and destruction]."
)
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "Status:"; Bold = $true },
    @{ Text = " "; Bold = $false },
    @{ Text = "None"; Bold = $false }
)
$idx = Add-Heading5 $d $idx "Project" "project-3"
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "Use of cryptographic key management for the Project system is in use for at the time of implementation for authentication. CivicActions utilizes customer agency supplied PIV credentials for access to customer instances of the Project. Access enforcement and authentication requirements for Project are described in AC-2 & IA-2. AWS platform does not utilize or manage cryptographic keys within the ACE boundary."; Bold = $false }
)

# ---------------------------------------------------------------------
# 6) Brand-new "SC-15: Collaborative Computing Devices" section, right
#    before "SC-20: Secure Name / Address Resolution Service
#    (Authoritative Source)".
# ---------------------------------------------------------------------

$idx = 52  # "SC-20: ..." Heading3 (shifted by +5 from step 5)
$idx = Add-Heading3 $d $idx "SC-15: Collaborative Computing Devices" "sc-15-collaborative-computing-devices"
$idx = Add-SourceCode $d $idx @(
    "The information system:",
    "  a.  Prohibits remote activation of collaborative computing devices with the",
    "following exceptions: [Assignment: organization-defined exceptions where remote activation is to be allowed]; and",
    "  b.  Provides an explicit indication of use to users physically present at the",
    "devices."
)
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "Status:"; Bold = $true },
    @{ Text = " "; Bold = $false },
    @{ Text = "None"; Bold = $false }
)
$idx = Add-Heading5 $d $idx "Project" "project-4"
$idx = Add-FirstParagraph $d $idx @(
    @{ Text = "This control is not applicable, as the Project system does"; Bold = $false },
    @{ Text = " "; Bold = $false },
    @{ Text = "employ any collaborative computing devices."; Bold = $false }
)

Write-Host "Done. Final paragraph count: " $d.Paragraphs.Count
